$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Taxas de mortalidade: corrected mortality-rate input for age 5 (M3) ---
$ws.Range("M3").Value2 = 0.00030030

# Column M (mortality-rate column) now needs the same "best fit" treatment
# previously applied only to column C.
$ws.Columns("M").ColumnWidth = $ws.Columns("C").ColumnWidth

# Cursor/selection left on H23 after the edit.
$ws.Range("H23").Select() | Out-Null

$excel.Calculate()
